$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Map of cell address -> new text value (all target cells are plain text,
# matching the original inline-string typed cells in the source sheet).
$updates = @{
    'D2' = '278.30'
    'E2' = '0.96%'
    'G2' = '19'
    'D3' = '27.46'
    'E3' = '0.50%'
    'G3' = '19'
    'E4' = '0.77%'
    'G4' = '19'
    'D5' = '0.06371'
    'E5' = '0.55%'
    'G5' = '19'
    'D6' = '7.031'
    'E6' = '1.24%'
    'G6' = '19'
    'D7' = '1.310'
    'E7' = '-2.70%'
    'G7' = '19'
    'D8' = '0.8941'
    'E8' = '1.90%'
    'G8' = '19'
    'E9' = '0.35%'
    'G9' = '19'
    'D10' = '0.05648'
    'E10' = '11.15%'
    'G10' = '19'
    'D11' = '0.07489'
    'E11' = '-1.25%'
    'G11' = '19'
    'D12' = '0.02911'
    'E12' = '-1.79%'
    'G12' = '19'
    'D13' = '0.08969'
    'E13' = '-0.62%'
    'G13' = '19'
    'D14' = '0.001571'
    'E14' = '0.27%'
    'G14' = '19'
    'D15' = '0.0006391'
    'E15' = '-0.36%'
    'G15' = '19'
    'D16' = '0.005976'
    'E16' = '1.63%'
    'G16' = '19'
    'E17' = '0.87%'
    'G17' = '19'
    'D18' = '3.324'
    'E18' = '0.83%'
    'G18' = '19'
    'D19' = '2.231'
    'E19' = '-1.79%'
    'G19' = '19'
    'G20' = '19'
    'E21' = '0.36%'
    'G21' = '19'
    'D22' = '3.899'
    'E22' = '-0.15%'
    'G22' = '19'
    'G23' = '19'
    'D24' = '0.04402'
    'E24' = '0.20%'
    'G24' = '19'
    'D25' = '0.001176'
    'E25' = '0.65%'
    'G25' = '19'
    'D26' = '0.004271'
    'E26' = '10.55%'
    'G26' = '19'
    'G27' = '19'
    'D28' = '0.0001180'
    'E28' = '-1.59%'
    'G28' = '19'
    'D29' = '0.0001653'
    'E29' = '-14.52%'
    'G29' = '19'
    'G30' = '19'
    'G31' = '19'
    'G32' = '19'
    'G33' = '19'
    'G34' = '19'
    'G35' = '19'
    'G36' = '19'
    'G37' = '19'
    'G38' = '19'
    'G39' = '19'
    'D40' = '0.04026'
    'E40' = '-3.61%'
    'G40' = '19'
    'D41' = '0.006724'
    'E41' = '-1.82%'
    'G41' = '19'
    'D42' = '0.1400'
    'E42' = '18.89%'
    'G42' = '19'
    'D43' = '0.002069'
    'E43' = '1.55%'
    'G43' = '19'
    'D44' = '0.01118'
    'E44' = '-3.07%'
    'G44' = '19'
    'D45' = '0.00005553'
    'E45' = '7.25%'
    'G45' = '19'
    'D46' = '1.628'
    'E46' = '9.53%'
    'G46' = '19'
    'D47' = '0.01849'
    'E47' = '-19.50%'
    'G47' = '19'
    'G48' = '19'
    'G49' = '19'
    'G50' = '19'
    'G51' = '19'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force a Text number format before assignment so Excel keeps the
    # literal string instead of re-interpreting "278.30" / "0.96%" / "19"
    # as a number/percentage (which would also drop trailing zeros).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Reset the style back to Normal so we do not leave a stray
    # number-format override on the cell (matches original formatting).
    $cell.Style = "Normal"
}
